$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell introduced in B1 (appended to shared strings as "victor")
$ws.Range("B1").Value = "victor"

# New questions appended after the existing A2:A11 block (rows 12-18)
$ws.Range("A12").Value = "Should the female start to flaps the mirror immediatly after a sensing a male or should it sense the beam first ?"

$ws.Range("A13").Value = "What is the exact light pattern for the male while searching ?"
$ws.Range("A13").Characters(19, 5).Font.Bold = $true
$ws.Range("A13").Characters(41, 5).Font.Bold = $true
$ws.Range("A13").Characters(51, 10).Font.Bold = $true

$ws.Range("A14").Value = "What is the exact female sound pattern to notify a male that she is interested ?"
$ws.Range("A14").Characters(19, 6).Font.Bold = $true
$ws.Range("A14").Characters(26, 5).Font.Bold = $true

$ws.Range("A15").Value = "What is the exact male encouragement sound pattern ?"
$ws.Range("A15").Characters(19, 4).Font.Bold = $true
$ws.Range("A15").Characters(38, 5).Font.Bold = $true

$ws.Range("A16").Value = "What is the exact climax light and sound pattern for the female ?"
$ws.Range("A16").Characters(26, 5).Font.Bold = $true
$ws.Range("A16").Characters(36, 5).Font.Bold = $true
$ws.Range("A16").Characters(58, 6).Font.Bold = $true

$ws.Range("A17").Value = "What is the exact climax light and sound pattern for the male ?"
$ws.Range("A17").Characters(26, 5).Font.Bold = $true
$ws.Range("A17").Characters(36, 5).Font.Bold = $true
$ws.Range("A17").Characters(58, 4).Font.Bold = $true

$ws.Range("A18").Value = "Does the mirror has a search state when the female doesn't remember the male ? (search state would be stop when the encourement signal is heard.)"
$ws.Range("A18").Characters(10, 6).Font.Bold = $true
$ws.Range("A18").Characters(52, 16).Font.Bold = $true

# Widen column A to fit the new, longer questions
$ws.Columns("A").ColumnWidth = 127.25

# Match the final selection recorded in the saved workbook
$ws.Range("B11").Select()
